$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" — append a new tracked file
# (b817ee4b-e125-40fb-920d-f425cf8f5474.md) as row 3 on each of the three
# sheets (Overview, zh-cn, de-de) and extend each sheet's table to match.
# ---------------------------------------------------------------------------

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/247c567f43a5d4253bfc0e26e2d7025e1deb49cd/e2e/b817ee4b-e125-40fb-920d-f425cf8f5474.md"

# --- Overview sheet --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "b817ee4b-e125-40fb-920d-f425cf8f5474.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-21 10:46:30"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $baseUrl, "", "", "e2e\b817ee4b-e125-40fb-920d-f425cf8f5474.md") | Out-Null

$tblOverview = $wsOverview.ListObjects.Item("Overview")
$tblOverview.Resize($wsOverview.Range("A1:G3"))

# --- zh-cn sheet -------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = "b817ee4b-e125-40fb-920d-f425cf8f5474.7041b905efc89bcd4bd9692aa55c16a1451c3e9f.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-21 10:46:26"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $baseUrl, "", "", "b817ee4b-e125-40fb-920d-f425cf8f5474.md") | Out-Null

$tblZh = $wsZh.ListObjects.Item("zh-cn")
$tblZh.Resize($wsZh.Range("A1:P3"))

# --- de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = "b817ee4b-e125-40fb-920d-f425cf8f5474.7041b905efc89bcd4bd9692aa55c16a1451c3e9f.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-21 10:46:30"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $baseUrl, "", "", "b817ee4b-e125-40fb-920d-f425cf8f5474.md") | Out-Null

$tblDe = $wsDe.ListObjects.Item("de-de")
$tblDe.Resize($wsDe.Range("A1:P3"))
